# Add run functions to run MLP, Logistic Regression and Random Forest with
# best parameters with kFold training/testing split. Add results on
# falsely-classified instances by these models.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------------
# 1. Row 3 (Majority Vote row) text is unchanged, only its rendered
#    height grew once the rows below it changed - match Excel's autofit.
# ---------------------------------------------------------------------
$ws.Rows(3).RowHeight = 102

# ---------------------------------------------------------------------
# 2. Row 4 ("...2 hidden layers...") gains a "(BEST, POST_THRESHOLD)"
#    marker and per-fold "N neurons per layer" best-parameter cells.
# ---------------------------------------------------------------------
$ws.Rows(4).RowHeight = 68
$ws.Range("A4").Value = "Multi-layer Perceptron`n    Photons Histogramizer: 6 buckets, filter late-arrival photons (BEST, POST_THRESHOLD)`n    MLP: 2 hidden layers, 8-40 neurons per layer"
$a4 = $ws.Range("A4")
$a4.Characters(90,22).Font.Underline = $true
$a4.Characters(90,22).Font.Name = "Calibri (Body)"
$a4.Characters(112,10).Font.Name = "Calibri"
$a4.Characters(122,15).Font.Bold = $true
$a4.Characters(122,15).Font.Name = "Calibri"
$a4.Characters(137,2).Font.Name = "Calibri"
$a4.Characters(139,5).Font.Underline = $true
$a4.Characters(139,5).Font.Name = "Calibri (Body)"
$a4.Characters(144,17).Font.Name = "Calibri (Body)"

$ws.Range("C4").Value = "39 neurons per layer"
$ws.Range("E4").Value = "15 neurons per layer"
$ws.Range("G4").Value = "10 neurons per layer"
$ws.Range("I4").Value = "24 neurons per layer"
$ws.Range("K4").Value = "36 neurons per layer"

# ---------------------------------------------------------------------
# 3. Row 5 ("...4 hidden layers...") gains the same marker plus its own
#    per-fold best-parameter cells.
# ---------------------------------------------------------------------
$ws.Rows(5).RowHeight = 68
$ws.Range("A5").Value = "Multi-layer Perceptron`n    Photons Histogramizer: 6 buckets, filter late-arrival photons (BEST, POST_THRESHOLD)`n    MLP: 4 hidden layers, 8-40 neurons per layer"
$a5 = $ws.Range("A5")
$a5.Characters(90,22).Font.Underline = $true
$a5.Characters(90,22).Font.Name = "Calibri (Body)"
$a5.Characters(112,10).Font.Name = "Calibri"
$a5.Characters(122,15).Font.Bold = $true
$a5.Characters(122,15).Font.Name = "Calibri"
$a5.Characters(137,2).Font.Name = "Calibri"
$a5.Characters(139,4).Font.Underline = $true
$a5.Characters(139,4).Font.Name = "Calibri (Body)"
$a5.Characters(143,18).Font.Name = "Calibri"

$ws.Range("C5").Value = "28 neurons per layer"
$ws.Range("E5").Value = "22 neurons per layer"
$ws.Range("G5").Value = "31 neurons per layer"
$ws.Range("I5").Value = "15 neurons per layer"
$ws.Range("K5").Value = "29 neurons per layer"

# ---------------------------------------------------------------------
# 4. Insert a brand-new row 6: the still-running 1-32 bucket MLP run.
# ---------------------------------------------------------------------
$ws.Rows(6).Insert()
$ws.Rows(6).RowHeight = 68
$ws.Range("A6").Value = "Multi-layer Perceptron`n    Photons Histogramizer: 1-32 buckets, filter late-arrival photons (POST_THRESHOLD)`n    MLP: 2 hidden layers, 8-40 neurons per layer"
$a6 = $ws.Range("A6")
$a6.Characters(51,4).Font.Underline = $true
$a6.Characters(51,4).Font.Name = "Calibri (Body)"
$a6.Characters(55,81).Font.Name = "Calibri"
$a6.Characters(136,4).Font.Underline = $true
$a6.Characters(136,4).Font.Name = "Calibri (Body)"
$a6.Characters(140,18).Font.Name = "Calibri"
$ws.Range("B6").Value = "(Still running)"

# ---------------------------------------------------------------------
# 5. Tighten up the blank gap below the table (old rows 7-15 shrink to
#    rows 7-10) so the trailing notes move from rows 16/18-21 to 11/13-16.
# ---------------------------------------------------------------------
$ws.Range("A7:A12").EntireRow.Delete()

# ---------------------------------------------------------------------
# 6. Text tweaks in the closing notes.
# ---------------------------------------------------------------------
$ws.Range("A15").Value = "Underscored parameters are those being tuned by a Grid Search. Each of the candidates is trained and we select the one that generates the best accuracy."
$ws.Range("A16").Value = "Threshold value for ""late-arrival phtoons"": have performed grid search on Max's previous empirical value (BEST) and Aaron's value (POST) shown by the histogram of all phtons' arrival times. Performance difference is insignificant."

# ---------------------------------------------------------------------
# 7. Column widths widen slightly now that longer "Accuracy" values sit
#    next to the new "Best Parameters" text.
# ---------------------------------------------------------------------
$ws.Columns(2).ColumnWidth = 11.1640625
$ws.Columns(4).ColumnWidth = 11.1640625
$ws.Columns(6).ColumnWidth = 11.1640625
$ws.Columns(8).ColumnWidth = 11.1640625
$ws.Columns(10).ColumnWidth = 11.1640625

# ---------------------------------------------------------------------
# 8. Selection moves to A6 (the newly-added row) as the last-touched cell.
# ---------------------------------------------------------------------
$ws.Range("A6").Select()
